# Citation-code normalization: replace ad-hoc reference placeholders with
# stable generated citation ids. Several placeholder codes are reused
# across different paragraphs but must map to *different* final ids
# depending on which paragraph (i.e. which source citation) they belong
# to, so each replacement is scoped to the specific paragraph's Range
# rather than done as a single document-wide Find/Replace.

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $range = $d.Paragraphs.Item($paraIndex).Range
    $range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2)
}

# Paragraph 3 - Sundiata / Sassouma Berete discussion: three citations,
# all normalized to the same new id (Ref-f457742).
Replace-InParagraph 3 "Ref-A1B2C3" "Ref-f457742"
Replace-InParagraph 3 "Ref-D4E5F6" "Ref-f457742"
Replace-InParagraph 3 "Ref-G7H8I9" "Ref-f457742"

# Paragraph 4 - Elizabeth / "A Question of Power" discussion.
Replace-InParagraph 4 "Ref-AB1CD2" "Ref-f494488"

# Paragraph 6 - Bofane / Congo Inc discussion.
Replace-InParagraph 6 "Ref-DJ7H2K" "Ref-f967437"

# Paragraph 7 - Dennis Brutus memoir discussion: two citations (same
# original placeholder codes as paragraph 3, but a different new id
# here), normalized to Ref-u247855.
Replace-InParagraph 7 "Ref-A1B2C3" "Ref-u247855"
Replace-InParagraph 7 "Ref-D4E5F6" "Ref-u247855"
